$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 117, pushing existing rows 117-165 down to 119-167
$ws.Range("A117:A118").EntireRow.Insert()

# Row 117
$ws.Cells.Item(117, 1).Value = 9
$ws.Cells.Item(117, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(117, 3).Value = "Metropolitana"
$ws.Cells.Item(117, 4).Value = 44876
$ws.Cells.Item(117, 5).Value = 13
$ws.Cells.Item(117, 6).Value = 300000000
$ws.Cells.Item(117, 7).Value = "Espárragos"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Banquete"
$ws.Cells.Item(117, 10).Value = 300
$ws.Cells.Item(117, 11).Value = 1500
$ws.Cells.Item(117, 12).Value = 1500
$ws.Cells.Item(117, 13).Value = 1500
$ws.Cells.Item(117, 14).Value = "$/kilo"
$ws.Cells.Item(117, 15).Value = "Provincia de Linares"
$ws.Cells.Item(117, 16).Value = 1500
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = "Hortaliza"

# Row 118
$ws.Cells.Item(118, 1).Value = 9
$ws.Cells.Item(118, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(118, 3).Value = "Metropolitana"
$ws.Cells.Item(118, 4).Value = 44876
$ws.Cells.Item(118, 5).Value = 13
$ws.Cells.Item(118, 6).Value = 300000000
$ws.Cells.Item(118, 7).Value = "Espárragos"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 550
$ws.Cells.Item(118, 11).Value = 1200
$ws.Cells.Item(118, 12).Value = 1300
$ws.Cells.Item(118, 13).Value = 1255
$ws.Cells.Item(118, 14).Value = "$/kilo"
$ws.Cells.Item(118, 15).Value = "Provincia de Linares"
$ws.Cells.Item(118, 16).Value = 1255
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = "Hortaliza"
